$d = $word.ActiveDocument

# 1) Insert " Хотилин М.И." right after "Преподаватель:" (end of that run),
#    preserving all existing runs untouched. A trailing placeholder
#    character is appended too, so the insertion point used for the
#    bookmark below is never the very last offset of the paragraph (a
#    collapsed range sitting immediately before the paragraph mark is
#    mis-resolved by bookmark placement) — it is removed afterwards.
$rng = $d.Content
$found = $rng.Find.Execute("Преподаватель:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(" Хотилин М.И.X")

# 2) Move the "_GoBack" bookmark from after "Лабораторная работа № 2" to the
#    end of the "Преподаватель: Хотилин М.И." paragraph (right after the
#    text we just inserted, covering the placeholder "X" for now). Adding a
#    bookmark with the same name removes/replaces the previous one at the
#    old location.
$target = $rng.Duplicate
$target.Collapse(0)
[void]$target.MoveStart(1, -1)
$d.Bookmarks.Add("_GoBack", $target)

# 3) Remove the placeholder character; the bookmark collapses back to sit
#    exactly between the inserted text and the paragraph mark.
$target.Delete()

# 4) The footer's cached PAGE field value flips from "8" to "2" (this doc is
#    now on page 2 instead of page 8). Range.Find inside the footer story
#    correctly *locates* the field result, but Range-offset based
#    Insert/Delete calls there are unreliable, so the field-result run's
#    text is updated through the Words collection instead (which maps back
#    to the right run correctly).
$footerRng = $d.Sections.Item(1).Footers.Item(1).Range
$pageWord = $footerRng.Paragraphs.Item(1).Range.Words.Item(1)
$pageWord.Text = "2"
